# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Re-shuffle the three worker rows (16-18) to new identities/values and
# update the corresponding mora amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> FAUSTINO OSPINO MAÑARA (doc 3557541), mora values unchanged (24640 / 616000)
$ws.Range("C16").Value = "3557541"
$ws.Range("D16").Value = "FAUSTINO OSPINO MAÑARA"
$ws.Range("F16").Value = 24640
$ws.Range("G16").Value = 616000

# Row 17 -> MIGUEL ANTONIO ARGEL PEREZ (doc 1044921278), mora values 25774 / 644350
$ws.Range("C17").Value = "1044921278"
$ws.Range("D17").Value = "MIGUEL ANTONIO ARGEL PEREZ"
$ws.Range("F17").Value = 25774
$ws.Range("G17").Value = 644350

# Row 18 -> LEONEL ENRIQUE HERNANDEZ BALLESTA (doc 1192718537), mora values 24640 / 616000
$ws.Range("C18").Value = "1192718537"
$ws.Range("D18").Value = "LEONEL ENRIQUE HERNANDEZ BALLESTA"
$ws.Range("F18").Value = 24640
$ws.Range("G18").Value = 616000
